# Generate Report for Handoff
# Adds a new "Ready for handoff" row (for file 96895750-271c-4152-856a-51bed20aeb0a...)
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$newFileName = "96895750-271c-4152-856a-51bed20aeb0aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newFilePath = "e2e\96895750-271c-4152-856a-51bed20aeb0aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newFileUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ecf4b5955f5ded33737af2906a1635fa97765a86/e2e/96895750-271c-4152-856a-51bed20aeb0aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

$status = "Ready for handoff"
$handoffDateTime = "2016-08-30 16:33:43"
$zhXliff = "96895750-271c-4152-856a-51bed20aeb0aoooooooooooooooooooooooooooooooooooooooo.462c7b5db04a75c8699ebe362e0b7c6cfebff1c8.zh-cn.xlf"
$zhXliffDateTime = "2016-08-30 16:33:38"
$deXliff = "96895750-271c-4152-856a-51bed20aeb0aoooooooooooooooooooooooooooooooooooooooo.462c7b5db04a75c8699ebe362e0b7c6cfebff1c8.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet (table "Overview") -> new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFileName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $handoffDateTime
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", $newFilePath) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet (table "zh_cn") -> new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $zhXliff
$wsZhCn.Range("H3").Value = $zhXliffDateTime
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newFileUrl, "", "", $newFileName) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet (table "de_de") -> new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $deXliff
$wsDeDe.Range("H3").Value = $handoffDateTime
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newFileUrl, "", "", $newFileName) | Out-Null
